$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Objetos")
Write-Host "E27:" $ws.Range("E27").Value()
Write-Host "E3:" $ws.Range("E3").Value()
Write-Host "F1:" $ws.Range("F1").Value()
Write-Host "Dim:" $ws.UsedRange.Address()
Write-Host "col C width:" $ws.Columns.Item(3).ColumnWidth
Write-Host "col D width:" $ws.Columns.Item(4).ColumnWidth
